# Update NATMI LR-pair TPM-derived metrics (Cntn1-Notch1) with newly
# recomputed values ("update scripts wuth new tpm").
# Columns: G=Ligand avg expr, H=Ligand total expr, M=Receptor avg expr,
# N=Receptor total expr, O/P=Receptor derived specificity (avg/total),
# Q/R=Edge avg/total expression weight, S/T=Edge derived specificity (avg/total)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.06624833333333334
$ws.Range("H2").Value = 0.198745
$ws.Range("M2").Value = 58.95713633333333
$ws.Range("N2").Value = 176.871409
$ws.Range("O2").Value = 0.4863146960083892
$ws.Range("P2").Value = 0.4863146960083893
$ws.Range("Q2").Value = 3.905812020189445
$ws.Range("R2").Value = 35.152308181705
$ws.Range("S2").Value = 0.4863146960083892
$ws.Range("T2").Value = 0.4863146960083893
$ws.Range("G3").Value = 0.06624833333333334
$ws.Range("H3").Value = 0.198745
$ws.Range("O3").Value = 0.07416766570679004
$ws.Range("P3").Value = 0.07416766570679005
$ws.Range("Q3").Value = 0.5956738766166667
$ws.Range("R3").Value = 5.36106488955
$ws.Range("S3").Value = 0.07416766570679004
$ws.Range("T3").Value = 0.07416766570679005
$ws.Range("G4").Value = 0.06624833333333334
$ws.Range("H4").Value = 0.198745
$ws.Range("M4").Value = 42.51661933333333
$ws.Range("N4").Value = 127.549858
$ws.Range("O4").Value = 0.3507032073181665
$ws.Range("P4").Value = 0.3507032073181665
$ws.Range("Q4").Value = 2.816655169801111
$ws.Range("R4").Value = 25.34989652821
$ws.Range("S4").Value = 0.3507032073181665
$ws.Range("T4").Value = 0.3507032073181665
$ws.Range("G5").Value = 0.06624833333333334
$ws.Range("H5").Value = 0.198745
$ws.Range("M5").Value = 10.76719366666667
$ws.Range("N5").Value = 32.301581
$ws.Range("O5").Value = 0.0888144309666542
$ws.Range("P5").Value = 0.08881443096665421
$ws.Range("Q5").Value = 0.7133086350938889
$ws.Range("R5").Value = 6.419777715845
$ws.Range("S5").Value = 0.0888144309666542
$ws.Range("T5").Value = 0.08881443096665421
